$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("F2").Value = 83
$ws.Range("H2").Value = 101
$ws.Range("J2").Value = 112
$ws.Range("F3").Value = 126
$ws.Range("G3").Value = 134
$ws.Range("J3").Value = 213
$ws.Range("B6").Value = 356
$ws.Range("E6").Value = 438
$ws.Range("H6").Value = 427
$ws.Range("J6").Value = 392
$ws.Range("B7").Value = 478
$ws.Range("E7").Value = 655
$ws.Range("F7").Value = 702
$ws.Range("G7").Value = 641
$ws.Range("H7").Value = 686
$ws.Range("J7").Value = 739

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("H6").Value = 36
$ws.Range("H7").Value = 47

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 7
$ws.Range("J7").Value = 43

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("B4").Value = 6
$ws.Range("B5").Value = 6

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J3").Value = 7
$ws.Range("J5").Value = 19

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("F2").Value = 6
$ws.Range("H6").Value = 25
$ws.Range("F7").Value = 52
$ws.Range("H7").Value = 44

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("H4").Value = 3
$ws.Range("H5").Value = 11

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("F21").Value = 11
$ws.Range("F28").Value = 52
$ws.Range("H28").Value = 44
$ws.Range("J29").Value = 12
$ws.Range("H32").Value = 47
$ws.Range("J36").Value = 43
$ws.Range("J41").Value = 19
$ws.Range("H48").Value = 5
$ws.Range("G53").Value = 80
$ws.Range("E76").Value = 19
$ws.Range("H82").Value = 11
$ws.Range("B88").Value = 6
$ws.Range("B98").Value = 478
$ws.Range("E98").Value = 655
$ws.Range("F98").Value = 702
$ws.Range("G98").Value = 641
$ws.Range("H98").Value = 686
$ws.Range("J98").Value = 739

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("G3").Value = 20
$ws.Range("G7").Value = 80

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("E5").Value = 9
$ws.Range("E6").Value = 19

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("J5").Value = 5
$ws.Range("J6").Value = 12

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("H2").Value = 1
$ws.Range("H6").Value = 5

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("F3").Value = 2
$ws.Range("F7").Value = 11
